$d = $word.ActiveDocument
$dash = [char]8211

# 1) "Alois Beran – [[PERSON_55]], [[PERSON_55]]" -> "[[PERSON_55]] – [[PERSON_55]], [[PERSON_55]]"
$find1 = "Alois Beran " + $dash + " [[PERSON_55]], [[PERSON_55]]"
$repl1 = "[[PERSON_55]] " + $dash + " [[PERSON_55]], [[PERSON_55]]"
$d.Content.Find.Execute($find1, $true, $false, $false, $false, $false, $true, 1, $false, $repl1, 2)

# 2) "[[PERSON_79]] – [[PERSON_79]], [[PERSON_80]]" -> "[[PERSON_79]] – [[PERSON_80]], [[PERSON_79]]"
$find2 = "[[PERSON_79]] " + $dash + " [[PERSON_79]], [[PERSON_80]]"
$repl2 = "[[PERSON_79]] " + $dash + " [[PERSON_80]], [[PERSON_79]]"
$d.Content.Find.Execute($find2, $true, $false, $false, $false, $false, $true, 1, $false, $repl2, 2)

# 3) "Max Kuchta – [[PERSON_82]], [[PERSON_82]]" -> "[[PERSON_82]] – [[PERSON_82]], [[PERSON_82]]"
$find3 = "Max Kuchta " + $dash + " [[PERSON_82]], [[PERSON_82]]"
$repl3 = "[[PERSON_82]] " + $dash + " [[PERSON_82]], [[PERSON_82]]"
$d.Content.Find.Execute($find3, $true, $false, $false, $false, $false, $true, 1, $false, $repl3, 2)

# 4) "Alex Kolísek – [[PERSON_84]], [[PERSON_84]]" -> "[[PERSON_84]] – [[PERSON_84]], [[PERSON_84]]"
$find4 = "Alex Kol" + [char]237 + "sek " + $dash + " [[PERSON_84]], [[PERSON_84]]"
$repl4 = "[[PERSON_84]] " + $dash + " [[PERSON_84]], [[PERSON_84]]"
$d.Content.Find.Execute($find4, $true, $false, $false, $false, $false, $true, 1, $false, $repl4, 2)
